# Applies the MEOW (Marine Ecoregions of the World) realm-name update:
# every existing realm label in column D (rows 284-512) gets an
# "_MEOW" suffix appended, including several rows whose original cell
# incorrectly referenced an unrelated country-name string rather than
# the intended realm name (e.g. row 309 previously showed "Bermuda",
# the country, instead of the MEOW realm "Bermuda"). The corrected,
# fully suffixed values are written out in row order below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$meowValues = @(
    'Adriatic Sea_MEOW',
    'Aegean Sea_MEOW',
    'Agulhas Bank_MEOW',
    'Alboran Sea_MEOW',
    'Aleutian Islands_MEOW',
    'Amazonia_MEOW',
    'Amsterdam-St Paul_MEOW',
    'Amundsen/Bellingshausen Sea_MEOW',
    'Andaman and Nicobar Islands_MEOW',
    'Andaman Sea Coral Coast_MEOW',
    'Angolan_MEOW',
    'Antarctic Peninsula_MEOW',
    'Arabian (Persian) Gulf_MEOW',
    'Arafura Sea_MEOW',
    'Araucanian_MEOW',
    'Arnhem Coast to Gulf of Carpenteria_MEOW',
    'Auckland Island_MEOW',
    'Azores Canaries Madeira_MEOW',
    'Baffin Bay - Davis Strait_MEOW',
    'Bahamian_MEOW',
    'Baltic Sea_MEOW',
    'Banda Sea_MEOW',
    'Bassian_MEOW',
    'Beaufort Sea - continental coast and shelf_MEOW',
    'Beaufort-Amundsen-Viscount Melville-Queen Maud_MEOW',
    'Bermuda_MEOW',
    'Bight of Sofala/Swamp Coast_MEOW',
    'Bismarck Sea_MEOW',
    'Black Sea_MEOW',
    'Bonaparte Coast_MEOW',
    'Bounty and Antipodes Islands_MEOW',
    'Bouvet Island_MEOW',
    'Campbell Island_MEOW',
    'Cape Howe_MEOW',
    'Cape Verde_MEOW',
    'Cargados Carajos/Tromelin Island_MEOW',
    'Carolinian_MEOW',
    'Celtic Seas_MEOW',
    'Central and Southern Great Barrier Reef_MEOW',
    'Central Chile_MEOW',
    'Central Kuroshio Current_MEOW',
    'Central New Zealand_MEOW',
    'Central Peru_MEOW',
    'Central Somali Coast_MEOW',
    'Chagos_MEOW',
    'Channels and Fjords of Southern Chile_MEOW',
    'Chatham Island_MEOW',
    'Chiapas-Nicaragua_MEOW',
    'Chiloense_MEOW',
    'Chukchi Sea_MEOW',
    'Clipperton_MEOW',
    'Cocos Islands_MEOW',
    'Cocos-Keeling/Christmas Island_MEOW',
    'Coral Sea_MEOW',
    'Cortezian_MEOW',
    'Crozet Islands_MEOW',
    'Delagoa_MEOW',
    'East African Coral Coast_MEOW',
    'East Antarctic Dronning Maud Land_MEOW',
    'East Antarctic Enderby Land_MEOW',
    'East Antarctic Wilkes Land_MEOW',
    'East Caroline Islands_MEOW',
    'East China Sea_MEOW',
    'East Greenland Shelf_MEOW',
    'East Siberian Sea_MEOW',
    'Easter Island_MEOW',
    'Eastern Bering Sea_MEOW',
    'Eastern Brazil_MEOW',
    'Eastern Caribbean_MEOW',
    'Eastern Galapagos Islands_MEOW',
    'Eastern India_MEOW',
    'Eastern Philippines_MEOW',
    'Exmouth to Broome_MEOW',
    'Faroe Plateau_MEOW',
    'Fernando de Naronha and Atoll das Rocas_MEOW',
    'Fiji Islands_MEOW',
    'Floridian_MEOW',
    'Gilbert/Ellis Islands_MEOW',
    'Great Australian Bight_MEOW',
    'Greater Antilles_MEOW',
    'Guayaquil_MEOW',
    'Guianan_MEOW',
    'Gulf of Aden_MEOW',
    'Gulf of Alaska_MEOW',
    'Gulf of Guinea Central_MEOW',
    'Gulf of Guinea Islands_MEOW',
    'Gulf of Guinea South_MEOW',
    'Gulf of Guinea Upwelling_MEOW',
    'Gulf of Guinea West_MEOW',
    'Gulf of Maine/Bay of Fundy_MEOW',
    'Gulf of Oman_MEOW',
    'Gulf of Papua_MEOW',
    'Gulf of St. Lawrence - Eastern Scotian Shelf_MEOW',
    'Gulf of Thailand_MEOW',
    'Gulf of Tonkin_MEOW',
    'Halmahera_MEOW',
    'Hawaii_MEOW',
    'Heard and Macdonald Islands_MEOW',
    'High Arctic Archipelago_MEOW',
    'Houtman_MEOW',
    'Hudson Complex_MEOW',
    'Humboldtian_MEOW',
    'Ionian Sea_MEOW',
    'Kamchatka Shelf and Coast_MEOW',
    'Kara Sea_MEOW',
    'Kerguelen Islands_MEOW',
    'Kermadec Island_MEOW',
    'Lancaster Sound_MEOW',
    'Laptev Sea_MEOW',
    'Leeuwin_MEOW',
    'Lesser Sunda_MEOW',
    'Levantine Sea_MEOW',
    'Line Islands_MEOW',
    'Lord Howe and Norfolk Islands_MEOW',
    'Macquarie Island_MEOW',
    'Magdalena Transition_MEOW',
    'Malacca Strait_MEOW',
    'Maldives_MEOW',
    'Malvinas/Falklands_MEOW',
    'Manning-Hawkesbury_MEOW',
    'Mariana Islands_MEOW',
    'Marquesas_MEOW',
    'Marshall Islands_MEOW',
    'Mascarene Islands_MEOW',
    'Mexican Tropical Pacific_MEOW',
    'Namaqua_MEOW',
    'Namib_MEOW',
    'Natal_MEOW',
    'New Caledonia_MEOW',
    'Nicoya_MEOW',
    'Ningaloo_MEOW',
    'North American Pacific Fijordland_MEOW',
    'North and East Barents Sea_MEOW',
    'North and East Iceland_MEOW',
    'North Greenland_MEOW',
    'North Patagonian Gulfs_MEOW',
    'North Sea_MEOW',
    'Northeast Sulawesi_MEOW',
    'Northeastern Brazil_MEOW',
    'Northeastern Honshu_MEOW',
    'Northeastern New Zealand_MEOW',
    'Northern and Central Red Sea_MEOW',
    'Northern Bay of Bengal_MEOW',
    'Northern California_MEOW',
    'Northern Galapagos Islands_MEOW',
    'Northern Grand Banks - Southern Labrador_MEOW',
    'Northern Gulf of Mexico_MEOW',
    'Northern Labrador_MEOW',
    'Northern Monsoon Current Coast_MEOW',
    'Northern Norway and Finnmark_MEOW',
    'Ogasawara Islands_MEOW',
    'Oregon, Washington, Vancouver Coast and Shelf_MEOW',
    'Oyashio Current_MEOW',
    'Palawan/North Borneo_MEOW',
    'Panama Bight_MEOW',
    'Papua_MEOW',
    'Patagonian Shelf_MEOW',
    'Peter the First Island_MEOW',
    'Phoenix/Tokelau/Northern Cook Islands_MEOW',
    'Puget Trough/Georgia Basin_MEOW',
    'Rapa-Pitcairn_MEOW',
    'Revillagigedos_MEOW',
    'Rio de la Plata_MEOW',
    'Rio Grande_MEOW',
    'Ross Sea_MEOW',
    'Saharan Upwelling_MEOW',
    'Sahelian Upwelling_MEOW',
    'Samoa Islands_MEOW',
    'Sao Pedro and Sao Paulo Islands_MEOW',
    'Scotian Shelf_MEOW',
    'Sea of Japan/East Sea_MEOW',
    'Sea of Okhotsk_MEOW',
    'Seychelles_MEOW',
    'Shark Bay_MEOW',
    'Snares Island_MEOW',
    'Society Islands_MEOW',
    'Solomon Archipelago_MEOW',
    'Solomon Sea_MEOW',
    'South and West Iceland_MEOW',
    'South Australian Gulfs_MEOW',
    'South European Atlantic Shelf_MEOW',
    'South Georgia_MEOW',
    'South India and Sri Lanka_MEOW',
    'South Kuroshio_MEOW',
    'South New Zealand_MEOW',
    'South Orkney Islands_MEOW',
    'South Sandwich Islands_MEOW',
    'South Shetland Islands_MEOW',
    'Southeast Madagascar_MEOW',
    'Southeast Papua New Guinea_MEOW',
    'Southeastern Brazil_MEOW',
    'Southern California Bight_MEOW',
    'Southern Caribbean_MEOW',
    'Southern China_MEOW',
    'Southern Cook/Austral Islands_MEOW',
    'Southern Grand Banks - South Newfoundland_MEOW',
    'Southern Gulf of Mexico_MEOW',
    'Southern Java_MEOW',
    'Southern Norway_MEOW',
    'Southern Red Sea_MEOW',
    'Southern Vietnam_MEOW',
    'Southwestern Caribbean_MEOW',
    'St. Helena and Ascension Islands_MEOW',
    'Sulawesi Sea/Makassar Strait_MEOW',
    'Sunda Shelf/Java Sea_MEOW',
    'Three Kings-North Cape_MEOW',
    'Tonga Islands_MEOW',
    'Torres Strait Northern Great Barrier Reef_MEOW',
    'Trindade and Martin Vaz Islands_MEOW',
    'Tristan Gough_MEOW',
    'Tuamotus_MEOW',
    'Tunisian Plateau/Gulf of Sidra_MEOW',
    'Tweed-Moreton_MEOW',
    'Uruguay-Buenos Aires Shelf_MEOW',
    'Vanuatu_MEOW',
    'Virginian_MEOW',
    'Weddell Sea_MEOW',
    'West Caroline Islands_MEOW',
    'West Greenland Shelf_MEOW',
    'Western and Northern Madagascar_MEOW',
    'Western Arabian Sea_MEOW',
    'Western Bassian_MEOW',
    'Western Caribbean_MEOW',
    'Western Galapagos Islands_MEOW',
    'Western India_MEOW',
    'Western Mediterranean_MEOW',
    'Western Sumatra_MEOW',
    'White Sea_MEOW',
    'Yellow Sea_MEOW'
)

$startRow = 284
for ($i = 0; $i -lt $meowValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $meowValues[$i]
}

"Updated " + $meowValues.Length + " cells in column D (rows " + $startRow + "-" + ($startRow + $meowValues.Length - 1) + ")"
